# Mobile App Documents are updated
# Updates the ECO FORM sheet with the new build/version/date info for the
# 01.02.02 release.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ECO FORM")

# Header "Date:" value -> Sept 21, 2015 (was Sept 02, 2015)
$ws.Range("H2").Value = "Sept 21, 2015"

# "Revision" value -> bump build version to 01.02.02
$ws.Range("H3").Value = "App Version - 01.02.01 , Build Version - 01.02.02"

# "Describe Change" release-note link text now points at the Active subfolder
$ws.Range("C9").Value = "RELEASENOTE_SRSMART_ANDROID_01.02.01.docx  - https://github.com/DelphianSystems/SecuRemote/tree/master/SR%20Smart%20App/Active`n"

# Testing date serial: Sept 2, 2015 -> Sept 21, 2015
$ws.Range("F15").Value = 42268

# Approvals / Upload-to-store dates
$ws.Range("G31").Value = "Date: 21/09/15"
$ws.Range("G33").Value = "Date: 21/09/15"

# Move the active selection as recorded in the saved workbook view
$ws.Range("K9").Select()
